$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the "Price" (D) and "Volume(1h)" (E) columns of the crypto table
# with the latest scraped figures. Both columns hold plain text (prices use
# "." as a thousands-style separator so they aren't real numbers, and the
# volume column keeps its leading/trailing padding spaces), so for any price
# that *would* parse as a genuine number we pin the cell to Text format
# first — otherwise Excel's normal type inference would silently turn
# "139.43" into the number 139.43 and drop things like trailing zeros
# (e.g. "0.110" -> 0.11).
$ws.Range("D2").Value = '57.582.09'
$ws.Range("E2").Value = '  +1.57%  '
$ws.Range("D3").Value = '3.012.68'
$ws.Range("E3").Value = '  +0.25%  '
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '509.12'
$ws.Range("E5").Value = '  -0.22%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '139.43'
$ws.Range("E6").Value = '  -0.07%  '
$ws.Range("E7").Value = '  +0.13%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.435'
$ws.Range("E8").Value = '  +0.44%  '
$ws.Range("E9").Value = '  +0.50%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.110'
$ws.Range("E10").Value = '  +1.18%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.365'
$ws.Range("E11").Value = '  +2.93%  '
$ws.Range("D12").Value = '3.533.39'
$ws.Range("E12").Value = '  +0.46%  '
$ws.Range("E13").Value = '  +0.73%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.54'
$ws.Range("E14").Value = '  +3.42%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000163'
$ws.Range("E15").Value = '  +4.65%  '
$ws.Range("D16").Value = '57.629.74'
$ws.Range("E16").Value = '  +1.54%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '6.25'
$ws.Range("E17").Value = '  +5.01%  '
$ws.Range("D18").Value = '3.018.54'
$ws.Range("E18").Value = '  +0.53%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.86'
$ws.Range("E19").Value = '  +2.89%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.95'
$ws.Range("E20").Value = '  +1.15%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '328.72'
$ws.Range("E21").Value = '  -0.24%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.998'
$ws.Range("E22").Value = '  -0.04%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.69'
$ws.Range("E23").Value = '  -1.51%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.499'
$ws.Range("E24").Value = '  +3.30%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '64.66'
$ws.Range("E25").Value = '  +2.81%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.167'
$ws.Range("E26").Value = '  -3.97%  '
$ws.Range("E27").Value = '  -0.27%  '
$ws.Range("D28").Value = '0.0₃0917'
$ws.Range("E28").Value = '  +0.38%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.77'
$ws.Range("E29").Value = '  +1.06%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.31'
$ws.Range("E30").Value = '  +3.35%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.81'
$ws.Range("E31").Value = '  +1.08%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.20'
$ws.Range("E32").Value = '  -5.00%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '20.60'
$ws.Range("E33").Value = '  -0.30%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.74'
$ws.Range("E34").Value = '  +3.38%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '153.99'
$ws.Range("E35").Value = '  -0.81%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.88'
$ws.Range("E36").Value = '  +3.57%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.27'
$ws.Range("E37").Value = '  +0.06%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '24.36'
$ws.Range("E38").Value = '  +0.24%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0676'
$ws.Range("E39").Value = '  -0.84%  '
$ws.Range("D40").Value = '3.049.04'
$ws.Range("E40").Value = '  +0.40%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '37.62'
$ws.Range("E41").Value = '  +1.78%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.85'
$ws.Range("E42").Value = '  +4.99%  '
$ws.Range("E43").Value = '  +0.03%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.650'
$ws.Range("E44").Value = '  +0.25%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.41'
$ws.Range("E45").Value = '  -0.37%  '
$ws.Range("D46").Value = '2.224.60'
$ws.Range("E46").Value = '  -1.90%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.986'
$ws.Range("E47").Value = '  -1.57%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '6.04'
$ws.Range("E48").Value = '  +3.71%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0238'
$ws.Range("E49").Value = '  -0.35%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '19.42'
$ws.Range("E50").Value = '  +0.11%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.86'
$ws.Range("E51").Value = '  -5.86%  '
